$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DatosGenerales")

# Set the "Historia" value (B6) to "23" and keep it as a text value
$ws.Range("B6").Value = "23"

# Update the active cell selection to B6 to match the saved view state
$ws.Range("B6").Select()
